# Apply the authors' update (2019-09-20) to row 6 of the bank-limit table:
#   - D6: single-transaction limit text changes from "5万" to "1万"
#   - G6: change-date moves from 2019-04-15 (43570) to 2019-09-20 (43728)
#   - active selection ends up on H5 (matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D6").Value = "1万"
$ws.Range("G6").Value = "2019-09-20"

$ws.Range("H5").Select()
